$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values replacing the previous Hello/World/123/123.456 data
$values = @("apple", "banana", "oranGe", "Grapes", "eggplant")

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $values[$i]
    $cell.Font.Bold = $false
}

# Remove the custom column width that was set on column A
$ws.Columns.Item(1).ColumnWidth = $ws.StandardWidth

# Remove any pictures/drawings (e.g. the logo image) from the sheet
for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
    $ws.Shapes.Item($i).Delete()
}
